# ------------------------------------------------------------------
# Add a new "2022-Q3" quarter sheet (with its fund-holding data) to the
# workbook, insert a matching summary row on the "总计" sheet, and keep
# the rest of the data (2021-Q4 / 2021-Q3 / 2020-Q4 sheets) untouched.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row right after the
#    header for the "2022-Q3" quarter, pushing the existing rows down.
# ------------------------------------------------------------------
$summary.Rows("2:2").Insert()
$summary.Range("A2:D2").ClearFormats()

# Copy the style of the index column (A) from the row below (still the
# plain "index" style used by every row in this column).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 19
$summary.Range("D2").Value = 2.12

# Re-number the plain row index stored in column A for the rows that
# got shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ------------------------------------------------------------------
# 2. Create the new "2022-Q3" worksheet right after "总计", matching the
#    layout used by the other quarterly sheets.
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Use the existing "2021-Q4" sheet (now shifted to index 3) purely as a
# formatting template for the header row / index column / page margins.
$refSheet = $wb.Worksheets.Item(3)
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$refSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundData = @(
    @("160106", "南方高增长混合（LOF）", "15.29", "91.51", "5.33", "0.8150", 5),
    @("160105", "南方积极配置混合（LOF）", "5.42", "91.29", "5.33", "0.2889", 5),
    @("001427", "招商丰泽灵活配置混合A", "9.64", "44.19", "1.93", "0.1861", 10),
    @("002819", "招商丰美灵活配置混合A", "7.73", "43.28", "1.92", "0.1484", 10),
    @("000314", "招商瑞丰灵活配置混合A", "5.78", "44.01", "1.93", "0.1116", 10),
    @("002389", "招商安德灵活配置混合A", "4.12", "45.53", "1.93", "0.0795", 10),
    @("000554", "南方中国梦灵活配置混合", "1.33", "90.51", "5.54", "0.0737", 4),
    @("002017", "招商瑞丰灵活配置混合C", "3.21", "44.01", "1.93", "0.0620", 10),
    @("015211", "招商安鼎平衡1年持有期混合A", "3.08", "44.14", "1.93", "0.0594", 10),
    @("014202", "天弘中证1000指数增强C", "3.69", "94.06", "1.56", "0.0576", 8),
    @("014201", "天弘中证1000指数增强A", "3.68", "94.06", "1.56", "0.0574", 8),
    @("002776", "招商安荣灵活配置混合A", "2.81", "43.84", "1.92", "0.0540", 10),
    @("015212", "招商安鼎平衡1年持有期混合C", "2.34", "44.14", "1.93", "0.0452", 10),
    @("002390", "招商安德灵活配置混合C", "1.50", "45.53", "1.93", "0.0290", 10),
    @("010434", "红土创新医疗保健股票", "0.44", "94.49", "4.91", "0.0216", 8),
    @("560006", "益民核心增长混合", "0.55", "77.02", "3.88", "0.0213", 1),
    @("001446", "招商丰泽灵活配置混合C", "0.28", "44.19", "1.93", "0.0054", 10),
    @("002820", "招商丰美灵活配置混合C", "0.21", "43.28", "1.92", "0.0040", 10),
    @("002777", "招商安荣灵活配置混合C", "0.19", "43.84", "1.92", "0.0036", 10),
)

$r = 2
foreach ($row in $fundData) {
    # Copy the index-column style down to every data row.
    $newSheet.Range("A2").Copy()
    $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $newSheet.Cells.Item($r, 1).Value = $r - 2

    # Fund code and the numeric-looking text columns must stay text
    # (otherwise Excel would strip leading zeros / turn them into
    # numbers), so force them in as text with a leading apostrophe and
    # then strip the style that operation adds back off again.
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 2).ClearFormats()

    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).ClearFormats()
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).ClearFormats()
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).ClearFormats()
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).ClearFormats()

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r++
}

Write-Output "workbook updated"
